$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a plain percentage string (e.g. "56%") need their number format
# forced to Text first, otherwise Excel auto-converts the input into a numeric
# percentage value (0.56) instead of keeping the literal text "56%".
$percentCells = @("H2", "H7", "H8", "H9", "H12", "H13", "H20", "H29", "H30", "H34", "H37", "H44", "H45", "H46")
foreach ($addr in $percentCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-17 18:48:28"
$ws.Range("H2").Value = "56%"
$ws.Range("E3").Value = "2026-02-17 18:48:30"
$ws.Range("O3").Value = "-4.2 °C"
$ws.Range("E4").Value = "2026-02-17 18:48:32"
$ws.Range("O4").Value = "9.5 °C"
$ws.Range("E5").Value = "2026-02-17 18:48:35"
$ws.Range("E6").Value = "2026-02-17 18:48:37"
$ws.Range("E7").Value = "2026-02-17 18:48:39"
$ws.Range("H7").Value = "60%"
$ws.Range("J7").Value = "1018.1 hPa"
$ws.Range("E8").Value = "2026-02-17 18:48:42"
$ws.Range("H8").Value = "72%"
$ws.Range("E9").Value = "2026-02-17 18:48:44"
$ws.Range("H9").Value = "55%"
$ws.Range("E10").Value = "2026-02-17 18:48:47"
$ws.Range("O10").Value = "10.7 °C"
$ws.Range("E11").Value = "2026-02-17 18:48:49"
$ws.Range("E12").Value = "2026-02-17 18:48:51"
$ws.Range("H12").Value = "57%"
$ws.Range("N12").Value = "9.7 °C 18:26 TU"
$ws.Range("O12").Value = "12.9 °C"
$ws.Range("E13").Value = "2026-02-17 18:48:54"
$ws.Range("H13").Value = "42%"
$ws.Range("E14").Value = "2026-02-17 18:48:56"
$ws.Range("E15").Value = "2026-02-17 18:48:59"
$ws.Range("N15").Value = "9.4 °C 18:26 TU"
$ws.Range("O15").Value = "12.5 °C"
$ws.Range("E16").Value = "2026-02-17 18:49:01"
$ws.Range("O16").Value = "-3.6 °C"
$ws.Range("E17").Value = "2026-02-17 18:49:03"
$ws.Range("E18").Value = "2026-02-17 18:49:06"
$ws.Range("E19").Value = "2026-02-17 18:49:08"
$ws.Range("E20").Value = "2026-02-17 18:49:10"
$ws.Range("H20").Value = "61%"
$ws.Range("O20").Value = "-2.0 °C"
$ws.Range("E21").Value = "2026-02-17 18:49:13"
$ws.Range("E22").Value = "2026-02-17 18:49:15"
$ws.Range("E23").Value = "2026-02-17 18:49:17"
$ws.Range("I23").Value = "2.6 mm"
$ws.Range("E24").Value = "2026-02-17 18:49:19"
$ws.Range("E25").Value = "2026-02-17 18:49:22"
$ws.Range("E26").Value = "2026-02-17 18:49:24"
$ws.Range("E27").Value = "2026-02-17 18:49:27"
$ws.Range("E28").Value = "2026-02-17 18:49:29"
$ws.Range("O28").Value = "8.9 °C"
$ws.Range("E29").Value = "2026-02-17 18:49:31"
$ws.Range("H29").Value = "63%"
$ws.Range("N29").Value = "8.7 °C 18:28 TU"
$ws.Range("E30").Value = "2026-02-17 18:49:33"
$ws.Range("H30").Value = "61%"
$ws.Range("O30").Value = "11.5 °C"
$ws.Range("E31").Value = "2026-02-17 18:49:36"
$ws.Range("J31").Value = "1018.3 hPa"
$ws.Range("E32").Value = "2026-02-17 18:49:38"
$ws.Range("O32").Value = "8.5 °C"
$ws.Range("E33").Value = "2026-02-17 18:49:41"
$ws.Range("J33").Value = "1017.2 hPa"
$ws.Range("E34").Value = "2026-02-17 18:49:43"
$ws.Range("H34").Value = "49%"
$ws.Range("O34").Value = "1.3 °C"
$ws.Range("E35").Value = "2026-02-17 18:49:45"
$ws.Range("E36").Value = "2026-02-17 18:49:48"
$ws.Range("E37").Value = "2026-02-17 18:49:50"
$ws.Range("H37").Value = "70%"
$ws.Range("J37").Value = "1018.8 hPa"
$ws.Range("E38").Value = "2026-02-17 18:49:52"
$ws.Range("E39").Value = "2026-02-17 18:49:55"
$ws.Range("K39").Value = "10.5 MJ/m2"
$ws.Range("O39").Value = "-2.5 °C"
$ws.Range("E40").Value = "2026-02-17 18:49:57"
$ws.Range("E41").Value = "2026-02-17 18:50:00"
$ws.Range("E42").Value = "2026-02-17 18:50:02"
$ws.Range("N42").Value = "10.5 °C 18:23 TU"
$ws.Range("O42").Value = "13.1 °C"
$ws.Range("E43").Value = "2026-02-17 18:50:04"
$ws.Range("O43").Value = "8.2 °C"
$ws.Range("E44").Value = "2026-02-17 18:50:07"
$ws.Range("H44").Value = "80%"
$ws.Range("E45").Value = "2026-02-17 18:50:09"
$ws.Range("H45").Value = "67%"
$ws.Range("I45").Value = "0.4 mm"
$ws.Range("J45").Value = "1021.6 hPa"
$ws.Range("E46").Value = "2026-02-17 18:50:12"
$ws.Range("H46").Value = "55%"
